$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.399168610572815
$ws.Range("B1").Value = 1.595760703086853
$ws.Range("C1").Value = 5.04764986038208
$ws.Range("D1").Value = 2.818417310714722
$ws.Range("E1").Value = 0.9069141149520874
